$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$styleBefore = $ws.Range("D11").Style
$ws.Range("D11:E24").Value = "'True"
$ws.Range("D11:E24").Style = $styleBefore
